$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 3925
$ws.Range("J58").Value = 5211.1113
$ws.Range("L58").Value = 15633.3339
$ws.Range("N58").Value = -15933.3339
$ws.Range("H69").Value = 8830.5
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 8830.5
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 26491.5
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -28239.5
$ws.Range("H72").Value = 8830.5
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 8830.5
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 79474.5
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -88210.5
$ws.Range("H92").Value = 3412.3333
$ws.Range("I92").Value = 2038.7646
$ws.Range("J92").Value = 9250
$ws.Range("K92").Value = 2038.7646
$ws.Range("L92").Value = 9250
$ws.Range("M92").Value = -790.7646
$ws.Range("N92").Value = -11746
$ws.Range("H93").Value = 70150.25
$ws.Range("J93").Value = 70150.25
$ws.Range("L93").Value = 70150.25
$ws.Range("N93").Value = -75142.25
$ws.Range("H137").Value = 4543.636
$ws.Range("I137").Value = 4543.636
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 13630.908
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -11080.908
$ws.Range("N137").ClearContents()
$ws.Range("H138").Value = 10795.24
$ws.Range("I138").Value = 10434.6
$ws.Range("J138").Value = 10814.221
$ws.Range("K138").Value = 31303.8
$ws.Range("L138").Value = 32442.663
$ws.Range("M138").Value = -26163.8
$ws.Range("N138").Value = -42722.663
$ws.Range("H141").Value = 3939.36
$ws.Range("I141").Value = 3939.36
$ws.Range("K141").Value = 11818.08
$ws.Range("M141").Value = -6638.08

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 1299.8
$ws.Range("J5").Value = 1750
$ws.Range("L5").Value = 1750
$ws.Range("N5").Value = -1974
$ws.Range("H32").Value = 19271.873
$ws.Range("I32").Value = 19271.873
$ws.Range("K32").Value = 19271.873
$ws.Range("M32").Value = -18984.873
$ws.Range("H74").Value = 2909.0588
$ws.Range("I74").Value = 2875.5
$ws.Range("J74").Value = 3446
$ws.Range("K74").Value = 2875.5
$ws.Range("L74").Value = 3446
$ws.Range("M74").Value = -2001.5
$ws.Range("N74").Value = -5194
$ws.Range("H77").Value = 2909.0588
$ws.Range("I77").Value = 2875.5
$ws.Range("J77").Value = 3446
$ws.Range("K77").Value = 14377.5
$ws.Range("L77").Value = 17230
$ws.Range("M77").Value = -10009.5
$ws.Range("N77").Value = -25966
$ws.Range("H96").Value = 29562
$ws.Range("J96").Value = 29562
$ws.Range("L96").Value = 29562
$ws.Range("N96").Value = -35054
$ws.Range("H122").Value = 6415.227
$ws.Range("I122").Value = 6415.227
$ws.Range("K122").Value = 19245.681
$ws.Range("M122").Value = -16795.681
$ws.Range("H132").Value = 12503115
$ws.Range("I132").Value = 3560.2856
$ws.Range("J132").Value = 100000000
$ws.Range("K132").Value = 10680.8568
$ws.Range("L132").Value = 300000000
$ws.Range("M132").Value = -8150.856800000001
$ws.Range("N132").Value = -300005060

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 1299.8
$ws.Range("J4").Value = 1750
$ws.Range("L4").Value = 1750
$ws.Range("N4").Value = -1980
$ws.Range("H82").Value = 30073.666
$ws.Range("I82").Value = 6088.4
$ws.Range("K82").Value = 6088.4
$ws.Range("M82").Value = -5705.4
$ws.Range("H85").Value = 30073.666
$ws.Range("I85").Value = 6088.4
$ws.Range("K85").Value = 6088.4
$ws.Range("M85").Value = -4762.4
$ws.Range("H99").Value = 1634.7
$ws.Range("I99").Value = 1705.8125
$ws.Range("J99").Value = 1350.25
$ws.Range("K99").Value = 1705.8125
$ws.Range("L99").Value = 1350.25
$ws.Range("M99").Value = -207.8125
$ws.Range("N99").Value = -4346.25
$ws.Range("H134").Value = 7145902
$ws.Range("I134").Value = 2820.5
$ws.Range("K134").Value = 8461.5
$ws.Range("M134").Value = -5926.5

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 40005104
$ws.Range("I31").Value = 76926696
$ws.Range("K31").Value = 76926696
$ws.Range("M31").Value = -76926401
$ws.Range("H34").Value = 40005104
$ws.Range("I34").Value = 76926696
$ws.Range("K34").Value = 76926696
$ws.Range("M34").Value = -76926494
$ws.Range("H59").Value = 95693.5
$ws.Range("J59").Value = 101665.664
$ws.Range("L59").Value = 101665.664
$ws.Range("N59").Value = -103955.664
$ws.Range("H86").Value = 6761.4
$ws.Range("I86").Value = 7033
$ws.Range("J86").Value = 6354
$ws.Range("K86").Value = 7033
$ws.Range("L86").Value = 6354
$ws.Range("M86").Value = -5910
$ws.Range("N86").Value = -8600
$ws.Range("H89").Value = 6761.4
$ws.Range("I89").Value = 7033
$ws.Range("J89").Value = 6354
$ws.Range("K89").Value = 35165
$ws.Range("L89").Value = 31770
$ws.Range("M89").Value = -29549
$ws.Range("N89").Value = -43002
$ws.Range("H122").Value = 4066.6667
$ws.Range("I122").Value = 4066.6667
$ws.Range("K122").Value = 12200.0001
$ws.Range("M122").Value = -9750.000100000001
$ws.Range("H132").Value = 2710.6
$ws.Range("I132").Value = 2710.6
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 8131.799999999999
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -5601.799999999999
$ws.Range("N132").ClearContents()
$ws.Range("H141").Value = 539910.6
$ws.Range("J141").Value = 577047.0600000001
$ws.Range("L141").Value = 577047.0600000001
$ws.Range("N141").Value = -587407.0600000001

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 80412
$ws.Range("J37").Value = 80412
$ws.Range("L37").Value = 241236
$ws.Range("N37").Value = -241460
$ws.Range("H137").Value = 11868.76
$ws.Range("I137").Value = 6001.5
$ws.Range("J137").Value = 13721.579
$ws.Range("K137").Value = 18004.5
$ws.Range("L137").Value = 41164.737
$ws.Range("M137").Value = -12904.5
$ws.Range("N137").Value = -51364.737
$ws.Range("H138").Value = 13548.177
$ws.Range("I138").Value = 13085.083
$ws.Range("K138").Value = 39255.249
$ws.Range("M138").Value = -34115.249
$ws.Range("H139").Value = 5468.1333
$ws.Range("I139").Value = 2566.1
$ws.Range("J139").Value = 11272.2
$ws.Range("K139").Value = 7698.299999999999
$ws.Range("L139").Value = 33816.60000000001
$ws.Range("M139").Value = -2558.299999999999
$ws.Range("N139").Value = -44096.60000000001
$ws.Range("H140").Value = 4387.737
$ws.Range("I140").Value = 2406.8333
$ws.Range("K140").Value = 7220.499899999999
$ws.Range("M140").Value = -2040.499899999999

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7529.8335
$ws.Range("I70").Value = 8007.8
$ws.Range("J70").Value = 5140
$ws.Range("K70").Value = 8007.8
$ws.Range("L70").Value = 5140
$ws.Range("M70").Value = -7737.8
$ws.Range("N70").Value = -5680
$ws.Range("H73").Value = 7529.8335
$ws.Range("I73").Value = 8007.8
$ws.Range("J73").Value = 5140
$ws.Range("K73").Value = 8007.8
$ws.Range("L73").Value = 5140
$ws.Range("M73").Value = -7071.8
$ws.Range("N73").Value = -7012
$ws.Range("H80").Value = 8834.166999999999
$ws.Range("I80").Value = 3005
$ws.Range("J80").Value = 10000
$ws.Range("K80").Value = 3005
$ws.Range("L80").Value = 10000
$ws.Range("M80").Value = -2007
$ws.Range("N80").Value = -11996
$ws.Range("H83").Value = 8834.166999999999
$ws.Range("I83").Value = 3005
$ws.Range("J83").Value = 10000
$ws.Range("K83").Value = 15025
$ws.Range("L83").Value = 50000
$ws.Range("M83").Value = -10033
$ws.Range("N83").Value = -59984
$ws.Range("H132").Value = 5618999
$ws.Range("I132").Value = 4099.579
$ws.Range("J132").Value = 12731204
$ws.Range("K132").Value = 12298.737
$ws.Range("L132").Value = 38193612
$ws.Range("M132").Value = -9768.736999999999
$ws.Range("N132").Value = -38198672

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 9489.079
$ws.Range("I40").Value = 6131.857
$ws.Range("K40").Value = 6131.857
$ws.Range("M40").Value = -5995.857

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 4730.0713
$ws.Range("I126").Value = 3672.3784
$ws.Range("J126").Value = 12557
$ws.Range("K126").Value = 11017.1352
$ws.Range("L126").Value = 37671
$ws.Range("M126").Value = -8547.135200000001
$ws.Range("N126").Value = -42611
$ws.Range("H132").Value = 669836.25
$ws.Range("I132").Value = 3396
$ws.Range("K132").Value = 10188
$ws.Range("M132").Value = -7658
